$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (id 1001) new master-data values -------------------------------
$ws.Range("B2").Value = "Vostro"
$ws.Range("C2").Value = "Dell"
$ws.Range("D2").Value = 3568
$ws.Range("E2").Value = "DKS"
$ws.Range("F2").Value = 1.454
$ws.Range("G2").Value = "To take enrollments"
$ws.Range("H2").Value = "eng"
# I2 (is_active), J2 (cr_by), K2 (cr_dtimes) stay the same values they had.

# --- Row 3 (id 1002) new master-data values (Arabic) -----------------------
$ws.Range("B3").Value = "ستر  "
$ws.Range("C3").Value = "دلّ  "
$ws.Range("D3").Value = 3568
$ws.Range("E3").Value = "DKS"
$ws.Range("F3").Value = 1.454
$ws.Range("G3").Value = "لأخذ التسجيلات"
$ws.Range("H3").Value = "ara"

# --- View: scroll right and select rows 4 downward (as done before saving) -
$ws.Range("C4").Select()
$ws.Rows("4:1048576").Select()

# --- Page setup tweaks (paper size = A4, portrait) -------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
